# This workbook holds weekly price-reporting rows (2-12) for
# "Vega Monumental Concepción - Albahaca". The underlying source data was
# refreshed, which reshuffles which week's figures (Fecha, Volumen, Precio
# mínimo/máximo/promedio ponderado, Origen, Precio $/Kg) land on which
# sheet row, while the descriptive columns (Mercado, Región, Categoría,
# Variedad, Calidad, Unidad de comercialización, Kg o Unidades,
# Clasificación) stay put.
#
# Build the new values per row (D, J, K, L, M, O, P) straight from the
# target state and write them back with the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per data row (row number => Fecha, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg)
$rows = @{
    2  = @{ D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana";   P = 463 }
    3  = @{ D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana";   P = 551 }
    4  = @{ D = 44659; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    5  = @{ D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana";   P = 484 }
    6  = @{ D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana";   P = 544 }
    7  = @{ D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana";   P = 548 }
    8  = @{ D = 44630; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    9  = @{ D = 44643; J = 90;  K = 2800; L = 3000; M = 2911; O = "Región Metropolitana";   P = 485 }
    10 = @{ D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana";   P = 622 }
    11 = @{ D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
    12 = @{ D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K: Precio mínimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Precio máximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio $/Kg
}
